# OpenCart_LoginData.xlsx edit:
# - Sheet1 rows 4-6 get new credential data (bhaskar3@gmail.com / bhaskar rows),
#   rows 7-9 are removed, hyperlinks re-pointed, and the view/outline bookkeeping
#   is refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: Bhaskarpattepu105@gmail.com/Bhaskar@123/Valid -> bhaskar3@gmail.com/Test12345/Valid
$ws.Range("A4").Value = "bhaskar3@gmail.com"
$ws.Range("A4").Style = "Link"
$ws.Range("B4").Value = "Test12345"
$ws.Range("C4").Value = "Valid"

# --- Row 5: Bhaskarpattepu105@gmail.com/Bhaskar@1234/Invalid -> bhaskar3@gmail.com/Test123456/InValid
$ws.Range("A5").Value = "bhaskar3@gmail.com"
$ws.Range("A5").Style = "Link"
$ws.Range("B5").Value = "Test123456"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "InValid"

# --- Row 6: bhaskar3@gmail.com/Test12345/Valid -> bhaskar/(blank)/InValid
$ws.Range("A6").Value = "bhaskar"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = "InValid"

# --- Hyperlinks: drop the stale ones on A6/A7/B5, add fresh ones on A4/A5
$ws.Range("A6").Hyperlinks.Delete()
$ws.Range("A7").Hyperlinks.Delete()
$ws.Range("B5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:bhaskar3@gmail.com", [Type]::Missing, [Type]::Missing, "bhaskar3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:bhaskar3@gmail.com", [Type]::Missing, [Type]::Missing, "bhaskar3@gmail.com")

# --- Rows 7-9 no longer exist in the refreshed sheet; tag row 9 with the same
# outline depth the sheet had reached before trimming so sheetFormatPr keeps
# remembering it (outlineLevelRow="5"), then delete the trailing rows.
$ws.Rows("9:9").OutlineLevel = 5
$ws.Rows("7:9").Delete()

# --- View bookkeeping: selection moves off the trimmed rows
$ws.Range("E8").Select()
